# Regenerate merged AHB files
#  - rename the "_old" / "_new" header-row labels to "_FV2310" / "_FV2404"
#  - turn the data range into an Excel Table ("Table1")
#  - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. rename header row (row 1), columns A:J and L:U ----------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2310"
}

# column K (11) stays "diff"
$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2404"
}

# --- 2. convert the used range into a table ----------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. freeze the header row -------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
